# Generate Report for Handoff
# Adds a new localization-status row (for 40dc2e36-8521-488d-aa7c-a48b230a2cec.md)
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/f4acab3dc01613719fe7d324877418ead3c20e4e/e2e/40dc2e36-8521-488d-aa7c-a48b230a2cec.md"

# Helper: write a literal "True"/"False" style text value into a cell without
# it being auto-coerced into a boolean cell. We do this by entering it as a
# quoted-string formula, then collapsing the formula down to its cached
# value in place (Copy + PasteSpecial values) so the final cell is a plain
# shared-string cell, matching how the sheet already stores these tokens.
function Set-TextValue($range, [string]$text) {
    $range.Formula = "=""" + $text + """"
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = "40dc2e36-8521-488d-aa7c-a48b230a2cec.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = "2016-08-13 00:49:03"
$wsOv.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $baseUrl, "", "", "e2e\40dc2e36-8521-488d-aa7c-a48b230a2cec.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
Set-TextValue $wsZh.Range("F3") "False"
$wsZh.Range("G3").Value = "40dc2e36-8521-488d-aa7c-a48b230a2cec.2b0db807c2695472885eea43fdfd75da5e4d7baf.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-13 00:48:53"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $wsZh.Range("M3") "True"
Set-TextValue $wsZh.Range("O3") "False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $baseUrl, "", "", "40dc2e36-8521-488d-aa7c-a48b230a2cec.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
Set-TextValue $wsDe.Range("F3") "False"
$wsDe.Range("G3").Value = "40dc2e36-8521-488d-aa7c-a48b230a2cec.2b0db807c2695472885eea43fdfd75da5e4d7baf.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-13 00:49:03"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $wsDe.Range("M3") "True"
Set-TextValue $wsDe.Range("O3") "False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $baseUrl, "", "", "40dc2e36-8521-488d-aa7c-a48b230a2cec.md") | Out-Null

Write-Output "Report rows added for 40dc2e36-8521-488d-aa7c-a48b230a2cec.md"
